# Update "想去人数" (column F) counts on all four worksheets
# (展览, 演出, 本地生活, 全部类型) to match the refreshed site export.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(3, 6).Value = 274
$ws.Cells.Item(4, 6).Value = 624
$ws.Cells.Item(5, 6).Value = 2761
$ws.Cells.Item(9, 6).Value = 6296
$ws.Cells.Item(11, 6).Value = 76
$ws.Cells.Item(13, 6).Value = 5003
$ws.Cells.Item(15, 6).Value = 543
$ws.Cells.Item(16, 6).Value = 2610
$ws.Cells.Item(17, 6).Value = 1347
$ws.Cells.Item(19, 6).Value = 1216
$ws.Cells.Item(20, 6).Value = 310
$ws.Cells.Item(22, 6).Value = 128
$ws.Cells.Item(23, 6).Value = 1061
$ws.Cells.Item(26, 6).Value = 531
$ws.Cells.Item(27, 6).Value = 1370
$ws.Cells.Item(28, 6).Value = 1034
$ws.Cells.Item(29, 6).Value = 2094
$ws.Cells.Item(30, 6).Value = 308
$ws.Cells.Item(31, 6).Value = 576
$ws.Cells.Item(32, 6).Value = 21
$ws.Cells.Item(33, 6).Value = 22
$ws.Cells.Item(34, 6).Value = 86
$ws.Cells.Item(36, 6).Value = 1490
$ws.Cells.Item(39, 6).Value = 113
$ws.Cells.Item(41, 6).Value = 15
$ws.Cells.Item(42, 6).Value = 292
$ws.Cells.Item(43, 6).Value = 2269
$ws.Cells.Item(44, 6).Value = 2548
$ws.Cells.Item(47, 6).Value = 271
$ws.Cells.Item(48, 6).Value = 103
$ws.Cells.Item(49, 6).Value = 92

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(8, 6).Value = 321
$ws.Cells.Item(10, 6).Value = 89
$ws.Cells.Item(15, 6).Value = 153
$ws.Cells.Item(16, 6).Value = 43
$ws.Cells.Item(19, 6).Value = 36
$ws.Cells.Item(23, 6).Value = 370
$ws.Cells.Item(24, 6).Value = 29
$ws.Cells.Item(32, 6).Value = 5

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(7, 6).Value = 568
$ws.Cells.Item(8, 6).Value = 1482
$ws.Cells.Item(9, 6).Value = 1811
$ws.Cells.Item(10, 6).Value = 2508
$ws.Cells.Item(11, 6).Value = 841
$ws.Cells.Item(12, 6).Value = 718
$ws.Cells.Item(13, 6).Value = 8

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(4, 6).Value = 274
$ws.Cells.Item(5, 6).Value = 624
$ws.Cells.Item(6, 6).Value = 568
$ws.Cells.Item(7, 6).Value = 2761
$ws.Cells.Item(9, 6).Value = 1482
$ws.Cells.Item(11, 6).Value = 2508
$ws.Cells.Item(12, 6).Value = 6296
$ws.Cells.Item(13, 6).Value = 841
$ws.Cells.Item(16, 6).Value = 5003
$ws.Cells.Item(17, 6).Value = 2610
$ws.Cells.Item(18, 6).Value = 1347
$ws.Cells.Item(20, 6).Value = 1216
$ws.Cells.Item(22, 6).Value = 128
$ws.Cells.Item(23, 6).Value = 321
$ws.Cells.Item(25, 6).Value = 89
$ws.Cells.Item(27, 6).Value = 1370
$ws.Cells.Item(28, 6).Value = 1034
$ws.Cells.Item(29, 6).Value = 2094
$ws.Cells.Item(30, 6).Value = 308
$ws.Cells.Item(31, 6).Value = 576
$ws.Cells.Item(32, 6).Value = 21
$ws.Cells.Item(33, 6).Value = 153
$ws.Cells.Item(34, 6).Value = 22
$ws.Cells.Item(36, 6).Value = 43
$ws.Cells.Item(39, 6).Value = 15
$ws.Cells.Item(42, 6).Value = 292
$ws.Cells.Item(43, 6).Value = 29
$ws.Cells.Item(44, 6).Value = 2269
$ws.Cells.Item(45, 6).Value = 2548
$ws.Cells.Item(47, 6).Value = 271
$ws.Cells.Item(48, 6).Value = 103
$ws.Cells.Item(49, 6).Value = 5

